$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Timestamp column (AD) for all data rows (2-52) to the new run timestamp.
# All AD cells currently share the same string value, matching the shared-string update in the diff.
$ws.Range("AD2:AD52").Value = "2024-09-20T16:21:37.102997"

# Row 2
$ws.Range("U2").Value = -3.5
$ws.Range("Y2").Value = 0.5

# Row 3
$ws.Range("U3").Value = 1
$ws.Range("Y3").Value = 1.5

# Row 4
$ws.Range("R4").Value = 51.5
$ws.Range("X4").Value = 0.0198019801980198
$ws.Range("AA4").Value = -0.0008924067822915333

# Row 6
$ws.Range("R6").Value = 51.5
$ws.Range("S6").Value = -110
$ws.Range("X6").Value = -0.07207207207207207
$ws.Range("AA6").Value = -0.1329805770813937

# Row 15
$ws.Range("U15").Value = -21.5
$ws.Range("Y15").Value = 44

# Row 16
$ws.Range("R16").Value = 48.5
$ws.Range("S16").Value = -110
$ws.Range("X16").Value = 0
$ws.Range("AA16").Value = 0.04865375020743345

# Row 18
$ws.Range("R18").Value = 60.5
$ws.Range("X18").Value = 0.03418803418803419
$ws.Range("AA18").Value = -0.04885430177258974

# Row 20
$ws.Range("R20").Value = 41.5
$ws.Range("S20").Value = -105
$ws.Range("X20").Value = -0.02352941176470588
$ws.Range("AA20").Value = -0.06912059935174901

# Row 28
$ws.Range("U28").Value = -1
$ws.Range("Y28").Value = -1.5

# Row 30
$ws.Range("U30").Value = -9
$ws.Range("Y30").Value = 0

# Row 33
$ws.Range("R33").Value = 53.5
$ws.Range("S33").Value = -105
$ws.Range("X33").Value = -0.01834862385321101
$ws.Range("AA33").Value = -0.09294822345008591

# Row 34
$ws.Range("R34").Value = 55.5
$ws.Range("X34").Value = -0.01769911504424779
$ws.Range("AA34").Value = 0.1935483870967742

# Row 45
$ws.Range("U45").Value = -27
$ws.Range("Y45").Value = -0.5

# Row 50
$ws.Range("U50").Value = -3.5
$ws.Range("Y50").Value = 9.5

# Row 25: clear Open/Current(T/U) and Move_s(Y) values entirely (cells removed in diff)
$ws.Range("T25").ClearContents()
$ws.Range("U25").ClearContents()
$ws.Range("Y25").ClearContents()

